$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44371
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("P2").Value = 438

# Row 3
$ws.Range("D3").Value = 44676
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = '$/caja 18 unidades'
$ws.Range("O3").Value = 'Región Metropolitana'
$ws.Range("P3").Value = 667
$ws.Range("Q3").Value = 18

# Row 4
$ws.Range("D4").Value = 44397
$ws.Range("J4").Value = 40
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 500

# Row 5
$ws.Range("D5").Value = 44355
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 500

# Row 6
$ws.Range("D6").Value = 44305
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("P6").Value = 438

# Row 7
$ws.Range("D7").Value = 44386
$ws.Range("J7").Value = 40
$ws.Range("O7").Value = 'Región del Maule'

# Row 8
$ws.Range("D8").Value = 44313
$ws.Range("J8").Value = 20

# Row 9
$ws.Range("D9").Value = 44392
$ws.Range("J9").Value = 95
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("N9").Value = '$/caja 16 unidades'
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 438
$ws.Range("Q9").Value = 16

# Row 10
$ws.Range("D10").Value = 44312
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("N10").Value = '$/caja 16 unidades'
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 438
$ws.Range("Q10").Value = 16

# Row 11
$ws.Range("D11").Value = 44467

# Row 12
$ws.Range("D12").Value = 44385
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("N12").Value = '$/caja 16 unidades'
$ws.Range("P12").Value = 438
$ws.Range("Q12").Value = 16

# Row 13
$ws.Range("D13").Value = 44420
$ws.Range("J13").Value = 45

# Row 14
$ws.Range("D14").Value = 44308
$ws.Range("J14").Value = 75
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 5000
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 312

# Row 15
$ws.Range("D15").Value = 44403
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = 5000
$ws.Range("P15").Value = 312

# Row 16
$ws.Range("D16").Value = 44348
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("O16").Value = 'Región del Maule'
$ws.Range("P16").Value = 438

# Row 17
$ws.Range("D17").Value = 44398
$ws.Range("J17").Value = 80
$ws.Range("O17").Value = 'Región Metropolitana'

# Row 18
$ws.Range("D18").Value = 44362
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 8000
$ws.Range("P18").Value = 500

# Row 19
$ws.Range("D19").Value = 44396
$ws.Range("J19").Value = 80
$ws.Range("O19").Value = 'Región Metropolitana'

# Row 20
$ws.Range("D20").Value = 44685
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("N20").Value = '$/caja 18 unidades'
$ws.Range("P20").Value = 667
$ws.Range("Q20").Value = 18

# Row 21
$ws.Range("D21").Value = 44399
$ws.Range("J21").Value = 80
$ws.Range("O21").Value = 'Región Metropolitana'

# Row 22
$ws.Range("D22").Value = 44389
$ws.Range("J22").Value = 55
$ws.Range("O22").Value = 'Región del Maule'

# Row 23
$ws.Range("D23").Value = 44301
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 12000
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = 44314
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = 5000
$ws.Range("P24").Value = 312

# Row 25
$ws.Range("D25").Value = 44354
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8500
$ws.Range("O25").Value = 'Región Metropolitana'
$ws.Range("P25").Value = 531

# Row 26
$ws.Range("D26").Value = 44354
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 9000
$ws.Range("P26").Value = 562

# Row 27
$ws.Range("D27").Value = 44694
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 12000
$ws.Range("N27").Value = '$/caja 18 unidades'
$ws.Range("P27").Value = 667
$ws.Range("Q27").Value = 18

# Row 28
$ws.Range("D28").Value = 44372
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 6400
$ws.Range("N28").Value = '$/caja 16 unidades'
$ws.Range("O28").Value = 'Región del Maule'
$ws.Range("P28").Value = 400
$ws.Range("Q28").Value = 16

# Row 29
$ws.Range("D29").Value = 44679
$ws.Range("J29").Value = 90
$ws.Range("K29").Value = 12000
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = 12000
$ws.Range("N29").Value = '$/caja 18 unidades'
$ws.Range("P29").Value = 667
$ws.Range("Q29").Value = 18

# Row 30
$ws.Range("D30").Value = 44369
$ws.Range("J30").Value = 60
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 7000
$ws.Range("M30").Value = 7000
$ws.Range("O30").Value = 'Región Metropolitana'
$ws.Range("P30").Value = 438

# Row 31
$ws.Range("D31").Value = 44315
$ws.Range("J31").Value = 40
